$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (E1, F1): copy style from D1, then set text values "256"/"512" ---
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("E1").Formula = "'256"
$ws.Range("F1").Formula = "'512"
$excel.CutCopyMode = $false

# --- Data rows 2-65: numeric values for E (x2) and F (x4) columns ---
$ws.Range("E2").Value = 344
$ws.Range("F2").Value = 688
$ws.Range("E3").Value = 14
$ws.Range("F3").Value = 14
$ws.Range("E4").Value = 344
$ws.Range("F4").Value = 688
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 15
$ws.Range("E6").Value = 354
$ws.Range("F6").Value = 709
$ws.Range("E7").Value = 14
$ws.Range("F7").Value = 14
$ws.Range("E8").Value = 380
$ws.Range("F8").Value = 761
$ws.Range("E9").Value = 14
$ws.Range("F9").Value = 15
$ws.Range("E10").Value = 344
$ws.Range("F10").Value = 688
$ws.Range("E11").Value = 14
$ws.Range("F11").Value = 14
$ws.Range("E12").Value = 344
$ws.Range("F12").Value = 688
$ws.Range("E13").Value = 17
$ws.Range("F13").Value = 17
$ws.Range("E14").Value = 358
$ws.Range("F14").Value = 720
$ws.Range("E15").Value = 14
$ws.Range("F15").Value = 14
$ws.Range("E16").Value = 368
$ws.Range("F16").Value = 732
$ws.Range("E17").Value = 14
$ws.Range("F17").Value = 14
$ws.Range("E18").Value = 344
$ws.Range("F18").Value = 688
$ws.Range("E19").Value = 344
$ws.Range("F19").Value = 688
$ws.Range("E20").Value = 344
$ws.Range("F20").Value = 688
$ws.Range("E21").Value = 344
$ws.Range("F21").Value = 688
$ws.Range("E22").Value = 367
$ws.Range("F22").Value = 734
$ws.Range("E23").Value = 367
$ws.Range("F23").Value = 734
$ws.Range("E24").Value = 351
$ws.Range("F24").Value = 705
$ws.Range("E25").Value = 351
$ws.Range("F25").Value = 705
$ws.Range("E26").Value = 344
$ws.Range("F26").Value = 688
$ws.Range("E27").Value = 344
$ws.Range("F27").Value = 688
$ws.Range("E28").Value = 344
$ws.Range("F28").Value = 688
$ws.Range("E29").Value = 344
$ws.Range("F29").Value = 688
$ws.Range("E30").Value = 354
$ws.Range("F30").Value = 795
$ws.Range("E31").Value = 354
$ws.Range("F31").Value = 795
$ws.Range("E32").Value = 349
$ws.Range("F32").Value = 703
$ws.Range("E33").Value = 349
$ws.Range("F33").Value = 703
$ws.Range("E34").Value = 26
$ws.Range("F34").Value = 26
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = 13
$ws.Range("E36").Value = 26
$ws.Range("F36").Value = 26
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = 13
$ws.Range("E38").Value = 24
$ws.Range("F38").Value = 24
$ws.Range("E39").Value = 12
$ws.Range("F39").Value = 12
$ws.Range("E40").Value = 26
$ws.Range("F40").Value = 26
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = 13
$ws.Range("E42").Value = 25
$ws.Range("F42").Value = 54
$ws.Range("E43").Value = 12
$ws.Range("F43").Value = 23
$ws.Range("E44").Value = 25
$ws.Range("F44").Value = 54
$ws.Range("E45").Value = 12
$ws.Range("F45").Value = 23
$ws.Range("E46").Value = 54
$ws.Range("F46").Value = 54
$ws.Range("E47").Value = 23
$ws.Range("F47").Value = 23
$ws.Range("E48").Value = 25
$ws.Range("F48").Value = 54
$ws.Range("E49").Value = 12
$ws.Range("F49").Value = 23
$ws.Range("E50").Value = 29
$ws.Range("F50").Value = 29
$ws.Range("E51").Value = 17
$ws.Range("F51").Value = 17
$ws.Range("E52").Value = 31
$ws.Range("F52").Value = 31
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = 13
$ws.Range("E54").Value = 29
$ws.Range("F54").Value = 29
$ws.Range("E55").Value = 17
$ws.Range("F55").Value = 17
$ws.Range("E56").Value = 31
$ws.Range("F56").Value = 31
$ws.Range("E57").Value = 13
$ws.Range("F57").Value = 13
$ws.Range("E58").Value = 31
$ws.Range("F58").Value = 31
$ws.Range("E59").Value = 12
$ws.Range("F59").Value = 12
$ws.Range("E60").Value = 31
$ws.Range("F60").Value = 31
$ws.Range("E61").Value = 12
$ws.Range("F61").Value = 12
$ws.Range("E62").Value = 31
$ws.Range("F62").Value = 31
$ws.Range("E63").Value = 12
$ws.Range("F63").Value = 12
$ws.Range("E64").Value = 31
$ws.Range("F64").Value = 31
$ws.Range("E65").Value = 12
$ws.Range("F65").Value = 12
